$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 27 (month 2026-02) stats
$ws.Range("B27").Value = 6536
$ws.Range("D27").Value = 6094181

# Recalculate dependent columns for row 27
$B27 = $ws.Range("B27").Value2
$C27 = $ws.Range("C27").Value2
$D27 = $ws.Range("D27").Value2

# 12 months earlier is row 15 (users_per_school / yoy comparisons)
$B15 = $ws.Range("B15").Value2
$C15 = $ws.Range("C15").Value2
$D15 = $ws.Range("D15").Value2

$ws.Range("E27").Value = $D27 / $B27
$ws.Range("F27").Value = (($B27 / $B15) - 1) * 100
$ws.Range("G27").Value = (($C27 / $C15) - 1) * 100
$ws.Range("H27").Value = (($D27 / $D15) - 1) * 100
